$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells stay as text (matches original inline-string
# formatting such as "210.60" / "1.560.08") instead of being auto-coerced
# into numeric values by the smart input parsing.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.305.31"
$ws.Range("E2").Value = "  -0.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.560.08"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.60"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.489"
$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.35"
$ws.Range("E8").Value = "  -4.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.55"
$ws.Range("E9").Value = "  -2.46%  "

$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0893"
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.554.73"
$ws.Range("E14").Value = "  -0.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.307.18"
$ws.Range("E15").Value = "  -0.78%  "

$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.512"
$ws.Range("E17").Value = "  -1.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.96"
$ws.Range("E18").Value = "  -1.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.31"
$ws.Range("E19").Value = "  -0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("E23").Value = "  +1.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.90"
$ws.Range("E24").Value = "  -2.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("E25").Value = "  -2.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.10"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.86"
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.33"
$ws.Range("E28").Value = "  -1.89%  "

$ws.Range("E29").Value = "  -0.52%  "

$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("E31").Value = "  +2.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.06"
$ws.Range("E32").Value = "  -4.27%  "

$ws.Range("E33").Value = "  -1.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("E34").Value = "  -1.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.375.70"
$ws.Range("E35").Value = "  -1.53%  "

$ws.Range("E36").Value = "  +1.98%  "

$ws.Range("E37").Value = "  -3.62%  "

$ws.Range("E38").Value = "  -0.50%  "

$ws.Range("E39").Value = "  +1.83%  "

$ws.Range("E40").Value = "  -2.05%  "

$ws.Range("E41").Value = "  -3.40%  "

$ws.Range("E42").Value = "  +2.15%  "

$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0472"
$ws.Range("E44").Value = "  -1.71%  "

$ws.Range("E45").Value = "  -1.27%  "

$ws.Range("E46").Value = "  -3.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.02"
$ws.Range("E47").Value = "  -1.30%  "

$ws.Range("E48").Value = "  -6.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.695.86"
$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.30"
$ws.Range("E50").Value = "  -1.17%  "

$ws.Range("E51").Value = "  -2.55%  "
